# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
#
# Sheets 1-3 and 5 use an "Ano" (Year) header prefix; sheet 4 uses an
# "Intervalo" (Interval) header prefix; sheet 6 only has a single year
# column (B1).

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo 2015"
$wsIntervalo.Range("C1").Value = "Intervalo 2015-2030"
$wsIntervalo.Range("D1").Value = "Intervalo 2031-2040"
$wsIntervalo.Range("E1").Value = "Intervalo 2041-2050"

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
